# Auto-generated Excel COM-interop script to apply scheduled data refresh
# to the Spriggan_Profits workbook (per diff / commit: "chore: update Sheets via scheduled runner")
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 14522.571
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("H19").Value = 1143.091
$ws.Range("I19").Value = 1116.5
$ws.Range("K19").Value = 1116.5
$ws.Range("M19").Value = -941.5
$ws.Range("H33").Value = 286.85715
$ws.Range("I33").Value = 268
$ws.Range("K33").Value = 268
$ws.Range("M33").Value = -39
$ws.Range("H62").Value = 1000
$ws.Range("I62").Value = 1000
$ws.Range("K62").Value = 1000
$ws.Range("M62").Value = -376
$ws.Range("H65").Value = 1000
$ws.Range("I65").Value = 1000
$ws.Range("K65").Value = 5000
$ws.Range("M65").Value = -1880
$ws.Range("H113").Value = 3992.5
$ws.Range("I113").Value = 3475
$ws.Range("J113").Value = 5027.5
$ws.Range("K113").Value = 3475
$ws.Range("L113").Value = 5027.5
$ws.Range("M113").Value = -221
$ws.Range("N113").Value = -11535.5
$ws.Range("N9").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2392.1738
$ws.Range("I2").Value = 2294.7144
$ws.Range("J2").Value = 2543.7778
$ws.Range("K2").Value = 2294.7144
$ws.Range("L2").Value = 2543.7778
$ws.Range("M2").Value = -2181.7144
$ws.Range("N2").Value = -2769.7778
$ws.Range("H116").Value = 2392.1738
$ws.Range("I116").Value = 2294.7144
$ws.Range("J116").Value = 2543.7778
$ws.Range("K116").Value = 2294.7144
$ws.Range("L116").Value = 2543.7778
$ws.Range("M116").Value = -0.7143999999998414
$ws.Range("N116").Value = -7131.7778

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2392.1738
$ws.Range("I3").Value = 2294.7144
$ws.Range("J3").Value = 2543.7778
$ws.Range("K3").Value = 2294.7144
$ws.Range("L3").Value = 2543.7778
$ws.Range("M3").Value = -2180.7144
$ws.Range("N3").Value = -2771.7778
$ws.Range("H51").Value = 77998
$ws.Range("J51").Value = 77998
$ws.Range("L51").Value = 77998
$ws.Range("N51").Value = -78980
$ws.Range("H99").Value = 2027.8334
$ws.Range("I99").Value = 2047
$ws.Range("J99").Value = 1989.5
$ws.Range("K99").Value = 2047
$ws.Range("L99").Value = 1989.5
$ws.Range("M99").Value = -549
$ws.Range("N99").Value = -4985.5
$ws.Range("H107").Value = 71774.13
$ws.Range("I107").Value = 5376.8
$ws.Range("J107").Value = 204568.8
$ws.Range("K107").Value = 5376.8
$ws.Range("L107").Value = 204568.8
$ws.Range("M107").Value = -3456.8
$ws.Range("N107").Value = -208408.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 26322566
$ws.Range("I58").Value = 50011200
$ws.Range("J58").Value = 1862.3334
$ws.Range("K58").Value = 50011200
$ws.Range("L58").Value = 1862.3334
$ws.Range("M58").Value = -50010997
$ws.Range("N58").Value = -2268.3334
$ws.Range("H99").Value = 17963.857
$ws.Range("I99").Value = 19053.54
$ws.Range("K99").Value = 19053.54
$ws.Range("M99").Value = -17555.54
$ws.Range("H107").Value = 551126.6
$ws.Range("I107").Value = 556089.4
$ws.Range("K107").Value = 556089.4
$ws.Range("M107").Value = -554169.4
$ws.Range("H126").Value = 17963.857
$ws.Range("I126").Value = 19053.54
$ws.Range("K126").Value = 57160.62
$ws.Range("M126").Value = -54690.62
$ws.Range("H136").Value = 26322566
$ws.Range("I136").Value = 50011200
$ws.Range("J136").Value = 1862.3334
$ws.Range("K136").Value = 150033600
$ws.Range("L136").Value = 5587.0002
$ws.Range("M136").Value = -150031050
$ws.Range("N136").Value = -10687.0002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 303
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 401.33334
$ws.Range("K2").Value = 48
$ws.Range("L2").Value = 2408.00004
$ws.Range("M2").Value = 65
$ws.Range("N2").Value = -2634.00004
$ws.Range("H23").Value = 3270
$ws.Range("J23").Value = 4850
$ws.Range("L23").Value = 14550
$ws.Range("N23").Value = -15020
$ws.Range("H34").Value = 949.5
$ws.Range("I34").Value = 949.5
$ws.Range("K34").Value = 2848.5
$ws.Range("M34").Value = -2764.5
$ws.Range("H55").Value = 500
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("H76").Value = 20000
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("H79").Value = 20000
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("M76").ClearContents()
$ws.Range("M79").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2211.8
$ws.Range("I80").Value = 2328.75
$ws.Range("J80").Value = 1744
$ws.Range("K80").Value = 2328.75
$ws.Range("L80").Value = 1744
$ws.Range("M80").Value = -1330.75
$ws.Range("N80").Value = -3740
$ws.Range("H83").Value = 2211.8
$ws.Range("I83").Value = 2328.75
$ws.Range("J83").Value = 1744
$ws.Range("K83").Value = 11643.75
$ws.Range("L83").Value = 8720
$ws.Range("M83").Value = -6651.75
$ws.Range("N83").Value = -18704
$ws.Range("H122").Value = 63270.6
$ws.Range("I122").Value = 72847.88
$ws.Range("J122").Value = 8999.333000000001
$ws.Range("K122").Value = 218543.64
$ws.Range("L122").Value = 26997.999
$ws.Range("M122").Value = -216093.64
$ws.Range("N122").Value = -31897.999
$ws.Range("H132").Value = 5440172
$ws.Range("I132").Value = 6583977
$ws.Range("K132").Value = 19751931
$ws.Range("M132").Value = -19749401

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1196.85
$ws.Range("I82").Value = 1250.25
$ws.Range("K82").Value = 1250.25
$ws.Range("M82").Value = -889.25
$ws.Range("H85").Value = 1196.85
$ws.Range("I85").Value = 1250.25
$ws.Range("K85").Value = 1250.25
$ws.Range("M85").Value = -2.25
$ws.Range("H114").Value = 40398
$ws.Range("J114").Value = 40398
$ws.Range("L114").Value = 40398
$ws.Range("N114").Value = -49076
$ws.Range("H132").Value = 15634309
$ws.Range("I132").Value = 15634309
$ws.Range("K132").Value = 46902927
$ws.Range("M132").Value = -46900397
$ws.Range("H136").Value = 1569.4
$ws.Range("I136").Value = 1299.6666
$ws.Range("J136").Value = 1974
$ws.Range("K136").Value = 3898.9998
$ws.Range("L136").Value = 5922
$ws.Range("M136").Value = -1348.9998
$ws.Range("N136").Value = -11022

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1702.4286
$ws.Range("I81").Value = 1152.8334
$ws.Range("K81").Value = 2305.6668
$ws.Range("M81").Value = -1244.6668
$ws.Range("H84").Value = 1702.4286
$ws.Range("I84").Value = 1152.8334
$ws.Range("K84").Value = 11528.334
$ws.Range("M84").Value = -6224.333999999999
$ws.Range("H120").Value = 60000
$ws.Range("J120").Value = 60000
$ws.Range("L120").Value = 60000
$ws.Range("N120").Value = -69676
$ws.Range("H121").Value = 99473.664
$ws.Range("J121").Value = 99473.664
$ws.Range("L121").Value = 99473.664
$ws.Range("N121").Value = -102967.664

Write-Host "Applied 179 cell updates/additions and 4 cell removals."